# Hortaliza, Vega Modelo de Temuco - Achicoria: add a new weekly price
# observation. A new record is inserted at row 127 (pushing the existing
# rows 127-152 down to 128-153), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127, shifting rows 127-152 down to 128-153.
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row 127 with the new observation.
$ws.Cells.Item(127, 1).Value = 10
$ws.Cells.Item(127, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(127, 3).Value = "La Araucanía"
$ws.Cells.Item(127, 4).Value = 45211
$ws.Cells.Item(127, 5).Value = 9
$ws.Cells.Item(127, 6).Value = 100112010
$ws.Cells.Item(127, 7).Value = "Achicoria"
$ws.Cells.Item(127, 8).Value = "Sin especificar"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 650
$ws.Cells.Item(127, 11).Value = 8000
$ws.Cells.Item(127, 12).Value = 10000
$ws.Cells.Item(127, 13).Value = 9077
$ws.Cells.Item(127, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(127, 15).Value = "Región Metropolitana"
$ws.Cells.Item(127, 16).Value = 504
$ws.Cells.Item(127, 17).Value = 18
$ws.Cells.Item(127, 18).Value = "Hortaliza"

# Note: Rows.Insert() already copies the row-above's per-cell formatting
# (including the "Fecha" column's date style) onto the new row 127, so no
# extra style assignment is needed here.
